$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (21st column) - shifts U,V -> V,W
$ws.Columns("U:U").Insert()

# New header cell at U1, matching the bold/centered style of the rest of row 1
$ws.Range("U1").Value = "MgCa Coretop modelled temperature"
$ws.Range("U1").Font.Bold = $true
$ws.Range("U1").HorizontalAlignment = -4108

# Update row 2 values per diff
$ws.Range("R2").Value = 27.01
$ws.Range("S2").Value = -0.01264004177517464
$ws.Range("T2").Value = -1.062940041775175
$ws.Range("U2").Value = 27.2019
$ws.Range("V2").Value = -0.2018999999999984
$ws.Range("W2").Value = -1.252199999999998
